$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 2-10: full cross-product of Sending cluster (A) x Target cluster (D)
# over {ECs, FAPs, sCs}, with fixed Ligand symbol (B) = Adam12, Receptor symbol (C) = Itga9.
$rows = @(
    @{A="ECs";  D="ECs";  E=2; F=0.6666666666666666; G=1.958141333333333;  H=5.874423999999999;  I=0.1445807708852573; J=0.1445807708852573; K=3; L=1; M=1.398034;           N=4.194102;  O=0.139066772576779;  P=0.139066772576779;  Q=2.737548160805333;  R=24.637933447248;    S=0.02010638118367546; T=0.02010638118367546},
    @{A="ECs";  D="FAPs"; E=2; F=0.6666666666666666; G=1.958141333333333;  H=5.874423999999999;  I=0.1445807708852573; J=0.1445807708852573; K=3; L=1; M=7.939250333333333;  N=23.817751; O=0.7897418235434783; P=0.7897418235434784; Q=15.54617423338044;  R=139.915568100424;   S=0.1141814816482449;  T=0.1141814816482449},
    @{A="ECs";  D="sCs";  E=2; F=0.6666666666666666; G=1.958141333333333;  H=5.874423999999999;  I=0.1445807708852573; J=0.1445807708852573; K=3; L=1; M=0.715685;            N=2.147055;  O=0.0711914038797426; P=0.0711914038797426; Q=1.401412380146666;  R=12.61271142132;     S=0.01029290805333688; T=0.01029290805333688},
    @{A="FAPs"; D="ECs";  E=3; F=1;                  G=5.833003000000001;  H=17.499009;          I=0.4306839633891008; J=0.4306839633891009; K=3; L=1; M=1.398034;           N=4.194102;  O=0.139066772576779;  P=0.139066772576779;  Q=8.154736516102;     R=73.39262864491801;  S=0.05989382878909788; T=0.0598938287890979},
    @{A="FAPs"; D="FAPs"; E=3; F=1;                  G=5.833003000000001;  H=17.499009;          I=0.4306839633891008; J=0.4306839633891009; K=3; L=1; M=7.939250333333333;  N=23.817751; O=0.7897418235434783; P=0.7897418235434784; Q=46.30967101208434;  R=416.787039108759;   S=0.3401291386178411;  T=0.3401291386178412},
    @{A="FAPs"; D="sCs";  E=3; F=1;                  G=5.833003000000001;  H=17.499009;          I=0.4306839633891008; J=0.4306839633891009; K=3; L=1; M=0.715685;           N=2.147055;  O=0.0711914038797426; P=0.0711914038797426; Q=4.174592752055;     R=37.571334768495;    S=0.03066099598216175; T=0.03066099598216175},
    @{A="sCs";  D="ECs";  E=3; F=1;                  G=5.752436333333333;  H=17.257309;          I=0.4247352657256419; J=0.4247352657256419; K=3; L=1; M=1.398034;           N=4.194102;  O=0.139066772576779;  P=0.139066772576779;  Q=8.042101576835332;  R=72.37891419151799;  S=0.05906656260400562; T=0.05906656260400563},
    @{A="sCs";  D="FAPs"; E=3; F=1;                  G=5.752436333333333;  H=17.257309;          I=0.4247352657256419; J=0.4247352657256419; K=3; L=1; M=7.939250333333333;  N=23.817751; O=0.7897418235434783; P=0.7897418235434784; Q=45.67003207689544;  R=411.0302886920589;  S=0.3354312032773922;  T=0.3354312032773923},
    @{A="sCs";  D="sCs";  E=3; F=1;                  G=5.752436333333333;  H=17.257309;          I=0.4247352657256419; J=0.4247352657256419; K=3; L=1; M=0.715685;           N=2.147055;  O=0.0711914038797426; P=0.0711914038797426; Q=4.116932397221666;  R=37.052391574995;    S=0.03023749984424396; T=0.03023749984424397}
)

$r = 2
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = "Adam12"
    $ws.Cells.Item($r, 3).Value = "Itga9"
    $ws.Cells.Item($r, 4).Value = $row.D
    $ws.Cells.Item($r, 5).Value = $row.E
    $ws.Cells.Item($r, 6).Value = $row.F
    $ws.Cells.Item($r, 7).Value = $row.G
    $ws.Cells.Item($r, 8).Value = $row.H
    $ws.Cells.Item($r, 9).Value = $row.I
    $ws.Cells.Item($r, 10).Value = $row.J
    $ws.Cells.Item($r, 11).Value = $row.K
    $ws.Cells.Item($r, 12).Value = $row.L
    $ws.Cells.Item($r, 13).Value = $row.M
    $ws.Cells.Item($r, 14).Value = $row.N
    $ws.Cells.Item($r, 15).Value = $row.O
    $ws.Cells.Item($r, 16).Value = $row.P
    $ws.Cells.Item($r, 17).Value = $row.Q
    $ws.Cells.Item($r, 18).Value = $row.R
    $ws.Cells.Item($r, 19).Value = $row.S
    $ws.Cells.Item($r, 20).Value = $row.T
    $r++
}
